# Adds the "LCD 8080 Pinout" section and a new E2 stepper group to the
# " Pin Function VET6" sheet, matching the commit "Added LCD 8080 Pinout".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(" Pin Function VET6")

# ---------------------------------------------------------------------
# 1) Shift the stepper signal labels in column E down one "axis group"
#    (X,Y,Z,E0,E1 -> Y,Z,E0,E1,E2) to make room for the new E2 stepper.
# ---------------------------------------------------------------------
$ws.Range("E2").Value  = "Y-STEP"
$ws.Range("E3").Value  = "Y-EN"
$ws.Range("E4").Value  = "Z-DIR"
$ws.Range("E5").Value  = "Z-STEP"
$ws.Range("E6").Value  = "Z-EN"
$ws.Range("E8").Value  = "E0-DIR"
$ws.Range("E9").Value  = "E0-STEP"
$ws.Range("E10").Value = "E0-EN"
$ws.Range("E11").Value = "E1-DIR"
$ws.Range("E12").Value = "E1-STEP"
$ws.Range("E13").Value = "E1-EN"
$ws.Range("E14").Value = "E2-DIR"
$ws.Range("E15").Value = "E2-STEP"
$ws.Range("E16").Value = "E2-EN"

# ---------------------------------------------------------------------
# 2) Row 40 moves from the SPI group into the LCD group.
# ---------------------------------------------------------------------
$ws.Range("A24:D24").Copy()
$ws.Range("A40:D40").PasteSpecial(-4122)
$ws.Range("D40").Value = "LCD"

# ---------------------------------------------------------------------
# 3) New LCD 8080 breakout pinout in column E, rows 24-40 (skipping the
#    SPI rows 35-39, which keep their own SPI-* signal names), plus the
#    merged title banner in G24:J24.
# ---------------------------------------------------------------------
$ws.Range("E24").Value = "LCD-D0"
$ws.Range("E25").Value = "LCD-D1"
$ws.Range("E26").Value = "LCD-D2"
$ws.Range("E27").Value = "LCD-D3"
$ws.Range("E28").Value = "LCD-D4"
$ws.Range("E29").Value = "LCD-D5"
$ws.Range("E30").Value = "LCD-D6"
$ws.Range("E31").Value = "LCD-D7"
$ws.Range("E32").Value = "LCD-CS"
$ws.Range("E33").Value = "LCD-CD"
$ws.Range("E34").Value = "LCD-WR"
$ws.Range("E35").Value = "SPI-CS1"
$ws.Range("E36").Value = "SPI-SCK"
$ws.Range("E37").Value = "SPI-MISO"
$ws.Range("E38").Value = "SPI-MOSI"
$ws.Range("E39").Value = "SPI-CS2"
$ws.Range("E40").Value = "LCD-RD"

$ws.Range("E24:E34").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E40").Copy()
$ws.Range("E40").PasteSpecial(-4122)

$ws.Range("A24:D24").Copy()
$ws.Range("G24").PasteSpecial(-4122)
$ws.Range("G24:J24").Merge()
$ws.Range("G24").Value = 'LCD Compatible with ADAFRUIT 2.8" 8080 Breakout'
$ws.Range("G24").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4) Rows 50/51 move from the ADD group into the EXTRA group (re-styled
#    to match the other EXTRA rows), and the EXTRA numbering downstream
#    (rows 58/59) shifts along with it.
# ---------------------------------------------------------------------
$ws.Range("A45:D45").Copy()
$ws.Range("A50:D50").PasteSpecial(-4122)
$ws.Range("A45:D45").Copy()
$ws.Range("A51:D51").PasteSpecial(-4122)
$ws.Range("D50").Value = "EXTRA"
$ws.Range("E50").Value = "EXTRA2"
$ws.Range("D51").Value = "EXTRA"
$ws.Range("E51").Value = "EXTRA3"

$ws.Range("E58").Value = "EXTRA4"
$ws.Range("E59").Value = "EXTRA5"

# ---------------------------------------------------------------------
# 5) Rows 70/71/72 move from the ADD/EXTRA group into the STEPPER group
#    (X-DIR/X-STEP/X-EN), and row 73 becomes Y-DIR.
# ---------------------------------------------------------------------
$ws.Range("A2:D2").Copy()
$ws.Range("A70:D70").PasteSpecial(-4122)
$ws.Range("A2:D2").Copy()
$ws.Range("A71:D71").PasteSpecial(-4122)
$ws.Range("A2:D2").Copy()
$ws.Range("A72:D72").PasteSpecial(-4122)

$ws.Range("D70").Value = "STEPPER"
$ws.Range("E70").Value = "X-DIR"
$ws.Range("D71").Value = "STEPPER"
$ws.Range("E71").Value = "X-STEP"
$ws.Range("D72").Value = "STEPPER"
$ws.Range("E72").Value = "X-EN"
$ws.Range("E73").Value = "Y-DIR"

# ---------------------------------------------------------------------
# 6) Sheet view bookkeeping to match the saved state after the edit.
# ---------------------------------------------------------------------
$ws.Range("E41").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
